# Scheduled market-data refresh: update the derived price/profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) for the
# affected leves across each world/job sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 3511.1
$ws.Cells.Item(62, 9).Value = 3184.1667
$ws.Cells.Item(62, 10).Value = 4001.5
$ws.Cells.Item(62, 11).Value = 3184.1667
$ws.Cells.Item(62, 12).Value = 4001.5
$ws.Cells.Item(62, 13).Value = -2560.1667
$ws.Cells.Item(62, 14).Value = -5249.5

$ws.Cells.Item(65, 8).Value = 3511.1
$ws.Cells.Item(65, 9).Value = 3184.1667
$ws.Cells.Item(65, 10).Value = 4001.5
$ws.Cells.Item(65, 11).Value = 15920.8335
$ws.Cells.Item(65, 12).Value = 20007.5
$ws.Cells.Item(65, 13).Value = -12800.8335
$ws.Cells.Item(65, 14).Value = -26247.5

$ws.Cells.Item(69, 8).Value = 171250.67
$ws.Cells.Item(69, 10).Value = 253751
$ws.Cells.Item(69, 12).Value = 761253
$ws.Cells.Item(69, 14).Value = -763001

$ws.Cells.Item(70, 8).Value = 1415.409
$ws.Cells.Item(70, 9).Value = 1188.7
$ws.Cells.Item(70, 10).Value = 1604.3334
$ws.Cells.Item(70, 11).Value = 3566.1
$ws.Cells.Item(70, 12).Value = 4813.0002
$ws.Cells.Item(70, 13).Value = -3296.1
$ws.Cells.Item(70, 14).Value = -5353.0002

$ws.Cells.Item(72, 8).Value = 171250.67
$ws.Cells.Item(72, 10).Value = 253751
$ws.Cells.Item(72, 12).Value = 2283759
$ws.Cells.Item(72, 14).Value = -2292495

$ws.Cells.Item(73, 8).Value = 1415.409
$ws.Cells.Item(73, 9).Value = 1188.7
$ws.Cells.Item(73, 10).Value = 1604.3334
$ws.Cells.Item(73, 11).Value = 3566.1
$ws.Cells.Item(73, 12).Value = 4813.0002
$ws.Cells.Item(73, 13).Value = -2630.1
$ws.Cells.Item(73, 14).Value = -6685.0002

$ws.Cells.Item(98, 8).Value = 2259.8572
$ws.Cells.Item(98, 9).Value = 2259.8572
$ws.Cells.Item(98, 11).Value = 2259.8572
$ws.Cells.Item(98, 13).Value = -761.8571999999999

$ws.Cells.Item(122, 8).Value = 2259.8572
$ws.Cells.Item(122, 9).Value = 2259.8572
$ws.Cells.Item(122, 11).Value = 6779.571599999999
$ws.Cells.Item(122, 13).Value = -4329.571599999999

$ws.Cells.Item(124, 8).Value = 79800
$ws.Cells.Item(124, 10).Value = 79800
$ws.Cells.Item(124, 12).Value = 79800
$ws.Cells.Item(124, 14).Value = -89620

$ws.Cells.Item(137, 8).Value = 2673.0757
$ws.Cells.Item(137, 9).Value = 1633.7805
$ws.Cells.Item(137, 10).Value = 4377.52
$ws.Cells.Item(137, 11).Value = 4901.3415
$ws.Cells.Item(137, 12).Value = 13132.56
$ws.Cells.Item(137, 13).Value = -2351.3415
$ws.Cells.Item(137, 14).Value = -18232.56

$ws.Cells.Item(138, 8).Value = 3956.948
$ws.Cells.Item(138, 9).Value = 2067.0386
$ws.Cells.Item(138, 10).Value = 4920.431
$ws.Cells.Item(138, 11).Value = 6201.1158
$ws.Cells.Item(138, 12).Value = 14761.293
$ws.Cells.Item(138, 13).Value = -1061.1158
$ws.Cells.Item(138, 14).Value = -25041.293

$ws.Cells.Item(141, 8).Value = 2906.6333
$ws.Cells.Item(141, 9).Value = 2350
$ws.Cells.Item(141, 10).Value = 3741.5833
$ws.Cells.Item(141, 11).Value = 7050
$ws.Cells.Item(141, 12).Value = 11224.7499
$ws.Cells.Item(141, 13).Value = -1870
$ws.Cells.Item(141, 14).Value = -21584.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 154039.55
$ws.Cells.Item(74, 9).Value = 163388
$ws.Cells.Item(74, 10).Value = 51206.5
$ws.Cells.Item(74, 11).Value = 163388
$ws.Cells.Item(74, 12).Value = 51206.5
$ws.Cells.Item(74, 13).Value = -162514
$ws.Cells.Item(74, 14).Value = -52954.5

$ws.Cells.Item(77, 8).Value = 154039.55
$ws.Cells.Item(77, 9).Value = 163388
$ws.Cells.Item(77, 10).Value = 51206.5
$ws.Cells.Item(77, 11).Value = 816940
$ws.Cells.Item(77, 12).Value = 256032.5
$ws.Cells.Item(77, 13).Value = -812572
$ws.Cells.Item(77, 14).Value = -264768.5

$ws.Cells.Item(122, 8).Value = 5001608.5
$ws.Cells.Item(122, 9).Value = 1609.7
$ws.Cells.Item(122, 10).Value = 25001602
$ws.Cells.Item(122, 11).Value = 4829.1
$ws.Cells.Item(122, 12).Value = 75004806
$ws.Cells.Item(122, 13).Value = -2379.1
$ws.Cells.Item(122, 14).Value = -75009706

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3348.4443
$ws.Cells.Item(31, 9).Value = 1767.8674
$ws.Cells.Item(31, 10).Value = 11547.6875
$ws.Cells.Item(31, 11).Value = 1767.8674
$ws.Cells.Item(31, 12).Value = 11547.6875
$ws.Cells.Item(31, 13).Value = -1472.8674
$ws.Cells.Item(31, 14).Value = -12137.6875

$ws.Cells.Item(34, 8).Value = 3348.4443
$ws.Cells.Item(34, 9).Value = 1767.8674
$ws.Cells.Item(34, 10).Value = 11547.6875
$ws.Cells.Item(34, 11).Value = 1767.8674
$ws.Cells.Item(34, 12).Value = 11547.6875
$ws.Cells.Item(34, 13).Value = -1565.8674
$ws.Cells.Item(34, 14).Value = -11951.6875

$ws.Cells.Item(58, 8).Value = 2276158
$ws.Cells.Item(58, 10).Value = 4456.6665
$ws.Cells.Item(58, 12).Value = 4456.6665
$ws.Cells.Item(58, 14).Value = -4862.6665

$ws.Cells.Item(94, 8).Value = 1361.375
$ws.Cells.Item(94, 9).Value = 1478.1666
$ws.Cells.Item(94, 10).Value = 1291.3
$ws.Cells.Item(94, 11).Value = 1478.1666
$ws.Cells.Item(94, 12).Value = 1291.3
$ws.Cells.Item(94, 13).Value = -1027.1666
$ws.Cells.Item(94, 14).Value = -2193.3

$ws.Cells.Item(99, 8).Value = 1660.8572
$ws.Cells.Item(99, 9).Value = 1702
$ws.Cells.Item(99, 10).Value = 1414
$ws.Cells.Item(99, 11).Value = 1702
$ws.Cells.Item(99, 12).Value = 1414
$ws.Cells.Item(99, 13).Value = -204
$ws.Cells.Item(99, 14).Value = -4410

$ws.Cells.Item(125, 8).Value = 60031.5
$ws.Cells.Item(125, 10).Value = 60031.5
$ws.Cells.Item(125, 12).Value = 60031.5
$ws.Cells.Item(125, 14).Value = -64951.5

$ws.Cells.Item(126, 8).Value = 1660.8572
$ws.Cells.Item(126, 9).Value = 1702
$ws.Cells.Item(126, 10).Value = 1414
$ws.Cells.Item(126, 11).Value = 5106
$ws.Cells.Item(126, 12).Value = 4242
$ws.Cells.Item(126, 13).Value = -2636
$ws.Cells.Item(126, 14).Value = -9182

$ws.Cells.Item(136, 8).Value = 2276158
$ws.Cells.Item(136, 10).Value = 4456.6665
$ws.Cells.Item(136, 12).Value = 13369.9995
$ws.Cells.Item(136, 14).Value = -18469.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 4561.7856
$ws.Cells.Item(34, 10).Value = 5084
$ws.Cells.Item(34, 12).Value = 15252
$ws.Cells.Item(34, 14).Value = -15420

$ws.Cells.Item(39, 8).Value = 9504.125
$ws.Cells.Item(39, 10).Value = 9504.125
$ws.Cells.Item(39, 12).Value = 28512.375
$ws.Cells.Item(39, 14).Value = -29100.375

$ws.Cells.Item(47, 8).Value = 301
$ws.Cells.Item(47, 9).Value = 301
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 903
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = -472
$ws.Cells.Item(47, 14).ClearContents()  # N47 removed (was -1612)

$ws.Cells.Item(55, 8).Value = 6500
$ws.Cells.Item(55, 10).Value = 10000
$ws.Cells.Item(55, 12).Value = 30000
$ws.Cells.Item(55, 14).Value = -30354

$ws.Cells.Item(113, 8).Value = 797.83673
$ws.Cells.Item(113, 9).Value = 797.2820400000001
$ws.Cells.Item(113, 10).Value = 800
$ws.Cells.Item(113, 11).Value = 2391.84612
$ws.Cells.Item(113, 12).Value = 2400
$ws.Cells.Item(113, 13).Value = -221.8461200000002
$ws.Cells.Item(113, 14).Value = -6740

$ws.Cells.Item(138, 8).Value = 7619.857
$ws.Cells.Item(138, 10).Value = 3866
$ws.Cells.Item(138, 12).Value = 11598
$ws.Cells.Item(138, 14).Value = -21878

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 11444
$ws.Cells.Item(40, 10).Value = 11444
$ws.Cells.Item(40, 12).Value = 11444
$ws.Cells.Item(40, 14).Value = -11746

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3634.975
$ws.Cells.Item(40, 9).Value = 3400.138
$ws.Cells.Item(40, 10).Value = 4254.091
$ws.Cells.Item(40, 11).Value = 3400.138
$ws.Cells.Item(40, 12).Value = 4254.091
$ws.Cells.Item(40, 13).Value = -3264.138
$ws.Cells.Item(40, 14).Value = -4526.091

$ws.Cells.Item(122, 8).Value = 4990.521
$ws.Cells.Item(122, 9).Value = 3998.3333
$ws.Cells.Item(122, 10).Value = 7967.0835
$ws.Cells.Item(122, 11).Value = 11994.9999
$ws.Cells.Item(122, 12).Value = 23901.2505
$ws.Cells.Item(122, 13).Value = -9544.999899999999
$ws.Cells.Item(122, 14).Value = -28801.2505

$ws.Cells.Item(136, 8).Value = 3772.7231
$ws.Cells.Item(136, 9).Value = 2502.2292
$ws.Cells.Item(136, 10).Value = 7360
$ws.Cells.Item(136, 11).Value = 7506.687600000001
$ws.Cells.Item(136, 12).Value = 22080
$ws.Cells.Item(136, 13).Value = -4956.687600000001
$ws.Cells.Item(136, 14).Value = -27180

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 5000
$ws.Cells.Item(52, 9).Value = 5000
$ws.Cells.Item(52, 11).Value = 5000
$ws.Cells.Item(52, 13).Value = -4774

$ws.Cells.Item(58, 8).Value = 16966.666
$ws.Cells.Item(58, 10).Value = 16966.666
$ws.Cells.Item(58, 12).Value = 16966.666
$ws.Cells.Item(58, 14).Value = -17582.666

$ws.Cells.Item(132, 8).Value = 2714.14
$ws.Cells.Item(132, 9).Value = 2263.342
$ws.Cells.Item(132, 10).Value = 4141.6665
$ws.Cells.Item(132, 11).Value = 6790.026
$ws.Cells.Item(132, 12).Value = 12424.9995
$ws.Cells.Item(132, 13).Value = -4260.026
$ws.Cells.Item(132, 14).Value = -17484.9995
